# TestSpreadsheet.xlsx — add the "Service Code" column (G) with its value "RSD"
# for the single data row, and leave the selection parked on the new header
# cell, matching the manual QA pass described in the commit
# ("Tested most workflows up to adding action data, fixed bugs").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell (row 1) and its data cell (row 2) in column G.
# These were previously empty (style-only) cells; they pick up two new
# shared-string entries ("Service Code" / "RSD") and keep their existing
# cell styles (s="1" / s="3") untouched.
$ws.Range("G1").Value = "Service Code"
$ws.Range("G2").Value = "RSD"

# Park the active selection on the newly added header cell.
$ws.Range("G1").Select()
